$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 191; existing rows 191:310 shift down to 192:311
$ws.Rows("191:191").Insert()

# Populate the newly inserted row 191 with the new weekly price record
$ws.Range("A191").Value = 4
$ws.Range("B191").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C191").Value = "Los Lagos"
$ws.Range("D191").Value = 44777
$ws.Range("E191").Value = 10
$ws.Range("F191").Value = 100112037
$ws.Range("G191").Value = "Cebollín"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 70
$ws.Range("K191").Value = 10000
$ws.Range("L191").Value = 11000
$ws.Range("M191").Value = 10500
$ws.Range("N191").Value = "$/paquete 36 unidades"
$ws.Range("O191").Value = "Región Metropolitana"
$ws.Range("P191").Value = 292
$ws.Range("Q191").Value = 36
$ws.Range("R191").Value = "Hortaliza"
